$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the date/owner cells for rows 27-29 (columns C, D, E) while keeping styles
$ws.Range("C27:E29").ClearContents()

# Update the view: scroll so row 16 is the top-left visible row, and change active selection to F24
$ws.Range("F24").Select()
$excel.ActiveWindow.ScrollRow = 16
